# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# coinranking.com snapshot. A few rows (40-44) also changed rank order,
# so Coin name + Link are rewritten for those rows too.
#
# Numeric-looking price strings (e.g. "3.00") are written with a leading
# "'" text-qualifier so Excel keeps them as literal text instead of
# coercing them to numbers (and dropping the trailing zero / decimals) -
# matching how these cells were already stored (t="inlineStr") before the
# edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.383.43'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '3.510.02'
$ws.Range("E3").Value = '  -3.59%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'" + '200.13'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = "'" + '553.42'
$ws.Range("E6").Value = '  -4.65%  '
$ws.Range("D7").Value = '3.500.81'
$ws.Range("E8").Value = '  -2.12%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  -3.57%  '
$ws.Range("D11").Value = "'" + '62.22'
$ws.Range("E11").Value = '  +10.52%  '
$ws.Range("E12").Value = '  -7.09%  '
$ws.Range("E13").Value = '  -7.44%  '
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").Value = '4.066.17'
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("D16").Value = '3.505.27'
$ws.Range("E16").Value = '  -3.82%  '
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = "'" + '18.49'
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").Value = '67.134.76'
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  -6.06%  '
$ws.Range("E21").Value = '  -5.50%  '
$ws.Range("D22").Value = "'" + '392.09'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = "'" + '12.29'
$ws.Range("E23").Value = '  -5.40%  '
$ws.Range("D24").Value = "'" + '4.01'
$ws.Range("E24").Value = '  -5.81%  '
$ws.Range("D25").Value = "'" + '83.13'
$ws.Range("E25").Value = '  -3.49%  '
$ws.Range("E26").Value = '  +2.76%  '
$ws.Range("D27").Value = "'" + '12.24'
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("E28").Value = '  -4.82%  '
$ws.Range("D29").Value = "'" + '8.85'
$ws.Range("E29").Value = '  -3.66%  '
$ws.Range("D30").Value = "'" + '31.00'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = "'" + '691.50'
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = "'" + '7.05'
$ws.Range("E32").Value = '  -12.96%  '
$ws.Range("D33").Value = "'" + '11.74'
$ws.Range("E33").Value = '  -4.22%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("E35").Value = '  -5.96%  '
$ws.Range("D36").Value = "'" + '38.78'
$ws.Range("E36").Value = '  -9.31%  '
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").Value = "'" + '0.397'
$ws.Range("E38").Value = '  -6.14%  '
$ws.Range("E39").Value = '  -5.12%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = "'" + '0.998'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").Value = "'" + '3.00'
$ws.Range("E41").Value = '  -4.51%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.063.79'
$ws.Range("E42").Value = '  -5.26%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0681'
$ws.Range("E43").Value = '  -13.61%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = "'" + '2.58'
$ws.Range("E44").Value = '  -9.82%  '
$ws.Range("E45").Value = '  +5.93%  '
$ws.Range("E46").Value = '  -4.45%  '
$ws.Range("E47").Value = '  -9.98%  '
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").Value = "'" + '138.14'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("E50").Value = '  -7.65%  '
$ws.Range("D51").Value = "'" + '2.88'
$ws.Range("E51").Value = '  -7.56%  '
